$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6428402662277222
$ws.Range("B1").Value = 0.7705091834068298
$ws.Range("C1").Value = 1.032574772834778
$ws.Range("D1").Value = 2.917504072189331
$ws.Range("E1").Value = 5.437038898468018
